$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two file-path strings to the new values
$ws.Range("B2").Value = "C:\Users\dpere\Documents\JTMT\Projects\תחזיות_דמוגרפיות\קבצי עבודה\142_מתחם_אנגל\בהת"
$ws.Range("B5").Value = "C:\Users\dpere\Documents\JTMT\forecast_by_version\V4"

# Row 2 height increases slightly
$ws.Rows.Item(2).RowHeight = 14.5

# Move the active selection from B4 to B5
$ws.Range("B5").Select() | Out-Null
